# Ran code for averaged intensities on spiral schemes.
#
# The "AlphaFiberF" sheet holds a table of averaging schemes (rows) vs HKL
# combinations (columns). We need to:
#   1. Insert 4 new rows right before the existing "NoRotation-tilt60deg"
#      row (row 10).
#   2. Re-home "Gaussian-Quadrature" (previously the very last scheme row)
#      into the first of those new rows, and add three new
#      "Spiral-..." scheme rows after it, all populated with 1s across the
#      HKL columns, matching the rest of the table.
#   3. Remove the old trailing "Gaussian-Quadrature" row (it has moved).
#   4. Renumber the column A index for every data row so it stays
#      sequential (0..17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the 4 new scheme rows just above "NoRotation-tilt60deg".
$ws.Rows("10:13").Insert()

# 2. The old "Gaussian-Quadrature" row got pushed down along with
#    everything else and now lives at row 20 - delete it, it has been
#    re-created above at row 10.
$ws.Rows("20:20").Delete()

# 3. Fill the newly inserted rows with the relocated / new schemes.
$newLabels = @("Gaussian-Quadrature", "Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space", "Spiral-90deg-10rot-3space")
for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $r = 10 + $i
    $ws.Cells.Item($r, 1).Value = 8 + $i
    $ws.Cells.Item($r, 2).Value = $newLabels[$i]
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# 4. Make column A's formatting (bold + border) consistent again for the
#    new rows - Insert() leaves them without the border style used by the
#    rest of the table.
$ws.Range("A9").Copy()
$ws.Range("A10:A13").PasteSpecial(-4122)

# 5. Renumber column A for the rows that shifted down so the index stays
#    sequential.
for ($r = 14; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
